$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: fill in previously-empty placeholder row ---
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 44745
$ws.Range("C42").Value = 0.61458333333333337
$ws.Range("D42").Value = 0.65625
$ws.Range("E42").Formula = "=D42-C42"
$ws.Range("F42").Value = "Code"
$ws.Range("G42").Value = "1. Block diagrams for Atrous convolutions, ASPP modules"

# --- Row 43: new row ---
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 44745
$ws.Range("C43").Value = 0.875
$ws.Range("D43").Value = 0.9375
$ws.Range("E43").Formula = "=D43-C43"
$ws.Range("F43").Value = "Code"
$ws.Range("G43").Value = "1. deeplabv3_plus_resnet50_starter nb completed`n2. deeplabv3_plus_resnet50_os8 nb completed"

$ws.Rows.Item(43).RowHeight = 30

# --- Update the selection to reflect where the author ended up ---
$ws.Range("G44").Select() | Out-Null

# --- Recalculate so E50 total and dependent formulas refresh ---
$excel.Calculate() | Out-Null
